$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("M2").Value = "nan"
$ws.Range("N2").Value = "nan"
$ws.Range("O2").Value = "nan"

$ws.Range("M4").Value = "nan"
